# Generate Report for Handback
#
# Refreshes the handback-status report after a new CI run:
#  - the existing tracked file (row 2 on every sheet) is re-stamped with
#    a freshly generated handback uuid/hash/timestamps
#    (d86cc8b9-... -> 1b31e8cc-...)
#  - a second handback file that was produced in the same run is
#    appended as a brand new row 3 (8c8f818e-...)

$wb = $excel.ActiveWorkbook

$uuid1 = "1b31e8cc-026b-48c4-b541-d3126b3474cc"
$uuid2 = "8c8f818e-736a-455b-af0a-3d9d4145e437"
$hash1 = "58f0a095076ec2beb27843689f2e5fcd17175051"
$hash2 = "b405e66c8165db295ed68ea946be9b34a58221d0"

$status = "Handed back: in sync with en-US"
$hyperColor = 15570276     # BGR encoding of RGB(100,149,237) = FF6495ED
$dateFormat = "yyyy-mm-dd HH:mm:ss"

function Style-AsLink($range) {
    $range.Font.Underline = $true
    $range.Font.Color = $hyperColor
}

function Style-AsDate($range) {
    $range.NumberFormat = $dateFormat
}

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ovTable = $ov.ListObjects.Item(1)

# refresh existing row 2 (was d86cc8b9-...)
$ov.Range("B2").Hyperlinks.Delete()
$ov.Range("A2").Value = "$uuid1.md"
$ov.Range("B2").Value = "e2e\$uuid1.md"
$ov.Range("C2").Value = ".md"
$ov.Range("E2").Value = $status
$ov.Range("F2").Value = $status
$ov.Range("G2").Value = "2016-09-03 21:04:32"
Style-AsDate $ov.Range("G2")
$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b5dbca17c3d468dc6f4e28dfc1e24f9b08ab383e/e2e/$uuid1.md", "", "", "e2e\$uuid1.md") | Out-Null
Style-AsLink $ov.Range("B2")

# add new row 3 (8c8f818e-...)
$ovTable.ListRows.Add() | Out-Null
$ov.Range("A3").Value = "$uuid2.md"
$ov.Range("B3").Value = "e2e\$uuid2.md"
$ov.Range("C3").Value = ".md"
$ov.Range("E3").Value = $status
$ov.Range("F3").Value = $status
$ov.Range("G3").Value = "2016-09-03 21:04:32"
Style-AsDate $ov.Range("G3")
$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b5dbca17c3d468dc6f4e28dfc1e24f9b08ab383e/e2e/$uuid2.md", "", "", "e2e\$uuid2.md") | Out-Null
Style-AsLink $ov.Range("B3")

# ---------------------------------------------------------------
# zh-cn / de-de detail sheets share the same column layout
# ---------------------------------------------------------------
function Update-LangSheet($ws, $langTag, $hoDate, $hbDate) {
    $lo = $ws.ListObjects.Item(1)

    # ---- refresh existing row 2 (uuid1) ----
    $ws.Range("A2").Hyperlinks.Delete()
    $ws.Range("I2").Hyperlinks.Delete()

    $ws.Range("A2").Value = "$uuid1.md"
    $ws.Range("B2").Value = ".md"
    $ws.Range("C2").Value = $status
    $ws.Range("D2").Value = "e2e"
    $ws.Range("E2").Value = "ht"
    $ws.Range("F2").Value = "False"
    $ws.Range("G2").Value = "$uuid1.$hash1.$langTag.xlf"
    $ws.Range("H2").Value = $hoDate
    Style-AsDate $ws.Range("H2")
    $ws.Range("I2").Value = "$uuid1.md"
    $ws.Range("J2").Value = "$uuid1.$hash1.$langTag.xlf"
    $ws.Range("K2").Value = $hbDate
    Style-AsDate $ws.Range("K2")
    $ws.Range("L2").Value = ""
    $ws.Range("M2").Value = "True"
    $ws.Range("N2").Value = ""
    $ws.Range("O2").Value = "False"
    $ws.Range("P2").Value = ""

    $ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b5dbca17c3d468dc6f4e28dfc1e24f9b08ab383e/e2e/$uuid1.md", "", "", "$uuid1.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-$($langTag.Replace('-',''))/blob/f44e8a81815bca2db23cfb72e40fc23f9ee9fa67/e2e/$uuid1.md", "", "", "$uuid1.md") | Out-Null
    Style-AsLink $ws.Range("A2")
    Style-AsLink $ws.Range("I2")

    # ---- add new row 3 (uuid2) ----
    $lo.ListRows.Add() | Out-Null
    $ws.Range("A3").Value = "$uuid2.md"
    $ws.Range("B3").Value = ".md"
    $ws.Range("C3").Value = $status
    $ws.Range("D3").Value = "e2e"
    $ws.Range("E3").Value = "ht"
    $ws.Range("F3").Value = "True"
    $ws.Range("G3").Value = "$uuid2.$hash2.$langTag.xlf"
    $ws.Range("H3").Value = $hoDate
    Style-AsDate $ws.Range("H3")
    $ws.Range("I3").Value = "$uuid2.md"
    $ws.Range("J3").Value = "$uuid2.$hash2.$langTag.xlf"
    $ws.Range("K3").Value = $hbDate
    Style-AsDate $ws.Range("K3")
    $ws.Range("L3").Value = ""
    $ws.Range("M3").Value = "True"
    $ws.Range("N3").Value = ""
    $ws.Range("O3").Value = "False"
    $ws.Range("P3").Value = ""

    $ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b5dbca17c3d468dc6f4e28dfc1e24f9b08ab383e/e2e/$uuid2.md", "", "", "$uuid2.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-$($langTag.Replace('-',''))/blob/f44e8a81815bca2db23cfb72e40fc23f9ee9fa67/e2e/$uuid2.md", "", "", "$uuid2.md") | Out-Null
    Style-AsLink $ws.Range("A3")
    Style-AsLink $ws.Range("I3")
}

$zh = $wb.Worksheets.Item("zh-cn")
Update-LangSheet $zh "zh-cn" "2016-09-03 21:04:27" "2016-09-03 21:04:45"

$de = $wb.Worksheets.Item("de-de")
Update-LangSheet $de "de-de" "2016-09-03 21:04:32" "2016-09-03 21:04:53"
